$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.400.12"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "2.007.08"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.639"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.16%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0741"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.103"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.897"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.63%  "

$ws.Range("E14").Value = "  +3.26%  "

$ws.Range("D15").Value = "2.300.88"
$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +15.49%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").Value = "2.052.19"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").Value = "36.314.55"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").Value = "0.0₃0860"
$ws.Range("E21").Value = "  +0.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +22.74%  "

$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  -1.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("E30").Value = "  -0.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.83%  "

$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +21.27%  "

$ws.Range("E34").Value = "  +4.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.08%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("E39").Value = "  +17.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +22.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.78%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.73%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0215"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.54%  "

$ws.Range("D48").Value = "1.434.12"
$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.84%  "
